$d = $word.ActiveDocument

# The document currently has three paragraphs:
#   1) "This is my great PhD thesis."
#   2) (empty)
#   3) (empty, carries Word's auto-managed "_GoBack" bookmark)
#
# The edit merges all three paragraphs into one: the existing sentence gets
# a new sentence appended to it (as its own separate run, not merged into
# the existing run's text), the "_GoBack" bookmark that used to sit alone
# in paragraph 3 now trails right after that text, and two brand-new empty
# paragraphs are left behind where the old paragraph breaks used to be.

$leftSingleQuote = [char]0x2018
$addedSentence = " But I " + $leftSingleQuote + "ate it."

# Grab the existing first paragraph's text and carry it forward verbatim.
# Paragraph.Range.Text includes the trailing paragraph-mark character
# (CR, chr 13), which is not part of the visible text, so trim it off
# before re-using the string. XML-escape before embedding below.
$existingText = $d.Paragraphs(1).Range.Text.TrimEnd([char]13)
$existingTextXml = $existingText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
$addedSentenceXml = $addedSentence.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$newParagraphXml = @"
<w:p $wordNs>
  <w:r><w:t>$existingTextXml</w:t></w:r>
  <w:r><w:t xml:space="preserve">$addedSentenceXml</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $wordNs/>
<w:p $wordNs/>
"@

# Replace the full span covering the three existing paragraphs (their
# paragraph marks included) with the new paragraph layout in one shot -
# this merges the paragraphs, appends the sentence as its own run, and
# carries the "_GoBack" bookmark over without disturbing it.
$target = $d.Range(0, $d.Content.End)
$target.InsertXML($newParagraphXml) | Out-Null
